# Scheduled runner update: refresh market/profit figures (columns H-N)
# on the leve-profit sheets (one sheet per crafting class).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2033.6842
$ws.Range("I40").Value = 2176.4546
$ws.Range("J40").Value = 1837.375
$ws.Range("K40").Value = 2176.4546
$ws.Range("L40").Value = 1837.375
$ws.Range("M40").Value = -2001.4546
$ws.Range("N40").Value = -2187.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3753.182
$ws.Range("J64").Value = 3100
$ws.Range("L64").Value = 3100
$ws.Range("N64").Value = -3596

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3753.182
$ws.Range("J67").Value = 3100
$ws.Range("L67").Value = 3100
$ws.Range("N67").Value = -4816

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 453.63635
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 453.63635
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1360.90905
$ws.Range("N80").Value = -3356.90905
$ws.Range("M80").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 453.63635
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 453.63635
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 4082.72715
$ws.Range("N83").Value = -14066.72715
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3578.5
$ws.Range("I86").Value = 6358.8
$ws.Range("J86").Value = 2651.7334
$ws.Range("K86").Value = 6358.8
$ws.Range("L86").Value = 2651.7334
$ws.Range("M86").Value = -5235.8
$ws.Range("N86").Value = -4897.7334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3578.5
$ws.Range("I89").Value = 6358.8
$ws.Range("J89").Value = 2651.7334
$ws.Range("K89").Value = 31794
$ws.Range("L89").Value = 13258.667
$ws.Range("M89").Value = -26178
$ws.Range("N89").Value = -24490.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 401.05
$ws.Range("I92").Value = 299.25
$ws.Range("J92").Value = 553.75
$ws.Range("K92").Value = 299.25
$ws.Range("L92").Value = 553.75
$ws.Range("M92").Value = 948.75
$ws.Range("N92").Value = -3049.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3182.647
$ws.Range("I106").Value = 3100.3845
$ws.Range("J106").Value = 3450
$ws.Range("K106").Value = 3100.3845
$ws.Range("L106").Value = 3450
$ws.Range("M106").Value = -2469.3845
$ws.Range("N106").Value = -4712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1945.4286
$ws.Range("I61").Value = 1663.0857
$ws.Range("K61").Value = 1663.0857
$ws.Range("M61").Value = -1451.0857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 57007.832
$ws.Range("I97").Value = 92011.82000000001
$ws.Range("J97").Value = 2001.5714
$ws.Range("K97").Value = 92011.82000000001
$ws.Range("L97").Value = 2001.5714
$ws.Range("M97").Value = -91515.82000000001
$ws.Range("N97").Value = -2993.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1945.4286
$ws.Range("I136").Value = 1663.0857
$ws.Range("K136").Value = 4989.257100000001
$ws.Range("M136").Value = -2439.257100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1893.95
$ws.Range("I86").Value = 1682.4615
$ws.Range("J86").Value = 2286.7144
$ws.Range("K86").Value = 1682.4615
$ws.Range("L86").Value = 2286.7144
$ws.Range("M86").Value = -559.4614999999999
$ws.Range("N86").Value = -4532.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1893.95
$ws.Range("I89").Value = 1682.4615
$ws.Range("J89").Value = 2286.7144
$ws.Range("K89").Value = 8412.307499999999
$ws.Range("L89").Value = 11433.572
$ws.Range("M89").Value = -2796.307499999999
$ws.Range("N89").Value = -22665.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4718908.5
$ws.Range("I134").Value = 7144080
$ws.Range("J134").Value = 3296.7222
$ws.Range("K134").Value = 21432240
$ws.Range("L134").Value = 9890.1666
$ws.Range("M134").Value = -21429705
$ws.Range("N134").Value = -14960.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1819.7241
$ws.Range("I132").Value = 1365.5238
$ws.Range("J132").Value = 3012
$ws.Range("K132").Value = 4096.5714
$ws.Range("L132").Value = 9036
$ws.Range("M132").Value = -1566.5714
$ws.Range("N132").Value = -14096

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1389.738
$ws.Range("I134").Value = 1250.5
$ws.Range("J134").Value = 1616
$ws.Range("K134").Value = 3751.5
$ws.Range("L134").Value = 4848
$ws.Range("M134").Value = -1216.5
$ws.Range("N134").Value = -9918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 7315
$ws.Range("J93").Value = 8478
$ws.Range("L93").Value = 25434
$ws.Range("N93").Value = -29178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 821.53125
$ws.Range("J122").Value = 823
$ws.Range("L122").Value = 7407
$ws.Range("N122").Value = -12307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 28990
$ws.Range("J26").Value = 28990
$ws.Range("L26").Value = 28990
$ws.Range("N26").Value = -29550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 28990
$ws.Range("J50").Value = 28990
$ws.Range("L50").Value = 28990
$ws.Range("N50").Value = -29986

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5202.6
$ws.Range("I70").Value = 4858.2856
$ws.Range("J70").Value = 6006
$ws.Range("K70").Value = 4858.2856
$ws.Range("L70").Value = 6006
$ws.Range("M70").Value = -4588.2856
$ws.Range("N70").Value = -6546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5202.6
$ws.Range("I73").Value = 4858.2856
$ws.Range("J73").Value = 6006
$ws.Range("K73").Value = 4858.2856
$ws.Range("L73").Value = 6006
$ws.Range("M73").Value = -3922.2856
$ws.Range("N73").Value = -7878

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 70554
$ws.Range("I80").Value = 3320.8
$ws.Range("J80").Value = 104170.6
$ws.Range("K80").Value = 3320.8
$ws.Range("L80").Value = 104170.6
$ws.Range("M80").Value = -2322.8
$ws.Range("N80").Value = -106166.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 70554
$ws.Range("I83").Value = 3320.8
$ws.Range("J83").Value = 104170.6
$ws.Range("K83").Value = 16604
$ws.Range("L83").Value = 520853
$ws.Range("M83").Value = -11612
$ws.Range("N83").Value = -530837

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1892.1765
$ws.Range("I122").Value = 1474.2727
$ws.Range("J122").Value = 2658.3333
$ws.Range("K122").Value = 4422.8181
$ws.Range("L122").Value = 7974.999899999999
$ws.Range("M122").Value = -1972.8181
$ws.Range("N122").Value = -12874.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 507.86365
$ws.Range("I22").Value = 246.66667
$ws.Range("J22").Value = 688.6923
$ws.Range("K22").Value = 246.66667
$ws.Range("L22").Value = 688.6923
$ws.Range("M22").Value = 48.33332999999999
$ws.Range("N22").Value = -1278.6923

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 507.86365
$ws.Range("I27").Value = 246.66667
$ws.Range("J27").Value = 688.6923
$ws.Range("K27").Value = 246.66667
$ws.Range("L27").Value = 688.6923
$ws.Range("M27").Value = -139.66667
$ws.Range("N27").Value = -902.6923

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2762.6
$ws.Range("I136").Value = 1367.5264
$ws.Range("J136").Value = 5172.273
$ws.Range("K136").Value = 4102.5792
$ws.Range("L136").Value = 15516.819
$ws.Range("M136").Value = -1552.5792
$ws.Range("N136").Value = -20616.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 38462210
$ws.Range("I107").Value = 45455100
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 136365300
$ws.Range("L107").Value = 3900
$ws.Range("M107").Value = -136363380
$ws.Range("N107").Value = -7740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3856.8044
$ws.Range("I136").Value = 1168.5938
$ws.Range("J136").Value = 10001.286
$ws.Range("K136").Value = 3505.7814
$ws.Range("L136").Value = 30003.858
$ws.Range("M136").Value = -955.7814000000003
$ws.Range("N136").Value = -35103.858
